# Remove the two "NA" rows (Durham School of Technology and Holton Career
# and Resource Center) from the high-school POC table. Deleting the entire
# rows shifts the remaining schools up, turning the original 14-row table
# (header + 13 schools/totals) into a 12-row table (header + 11
# schools/totals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 first (Holton Career and Resource Center) so row 8's index is not
# disturbed before it is removed too.
$ws.Rows(10).EntireRow.Delete()
$ws.Rows(8).EntireRow.Delete()

# Match the saved cursor position from the source file.
$ws.Range("C11").Select()
